$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (the "Чиизу Рамен" item) was incorrectly reusing the GUID string from
# row 7 ("Дайзу рамен") in column C. Give it its own fresh identifier so the
# update/import logic (which matches on this id) works correctly again.
$ws.Cells.Item(9, 3).Value2 = "45642981-18cb-4716-9cce-4ec65f149555"

# Widen column C so the (now-visible) GUID column reads comfortably.
$ws.Columns.Item(3).ColumnWidth = 38.33
